# fix(gui) step 1 and 2
# - bump the price list's date (A1) by one day
# - apply the new (step 2) prices to the "step 1/2" price columns

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 holds the list date as an Excel serial date: 45308 (2024-01-17) -> 45309 (2024-01-18)
$ws.Range("A1").Value = 45309

# Updated prices (15% increase) in column D
$ws.Range("D22").Value = 364.361
$ws.Range("D23").Value = 462.221
$ws.Range("D34").Value = 360.119
$ws.Range("D35").Value = 498.753
$ws.Range("D45").Value = 447.541
$ws.Range("D46").Value = 495.818
